$d = $word.ActiveDocument

# Map of exact current paragraph text -> new paragraph text.
# Only the option-list runs ("a) ", "b) ", "c) ", "d) " prefixes) are
# stripped; the "Resposta: ..." answer lines are left untouched.
$map = @{
    "a) Kelly "      = "Kelly "
    "b) Yorkie "     = "Yorkie "
    "c) Susan "      = "Susan "
    "d) Martha"      = "Martha"
    "a) Taffy's "    = "Taffy's "
    "b) Tucker's "   = "Tucker's "
    "c) Tucker's 2 " = "Tucker's 2 "
    "d) Taffy's 2"   = "Taffy's 2"
}

foreach ($p in $d.Paragraphs) {
    $txt = $p.Range.Text.TrimEnd([char]13)
    if ($map.ContainsKey($txt)) {
        $p.Range.Text = $map[$txt]
    }
}
